# Update LeveProfits pricing data across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block rewrites the H:N (currentAveragePrice..LeveProfitHQ) span for a specific row using a
# single array assignment; $null entries clear/omit a cell (matching rows where a column has no value).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$r6ALC = New-Object 'object[,]' 1,7
$r6ALC[0,0] = 281.55554
$r6ALC[0,1] = 191.875
$r6ALC[0,2] = 999
$r6ALC[0,3] = 575.625
$r6ALC[0,4] = 2997
$r6ALC[0,5] = -463.625
$r6ALC[0,6] = -3221
$ws.Range("H6:N6").Value = $r6ALC

$r8ALC = New-Object 'object[,]' 1,6
$r8ALC[0,0] = 244.91667
$r8ALC[0,1] = 226.55556
$r8ALC[0,2] = 300
$r8ALC[0,3] = 679.66668
$r8ALC[0,4] = 900
$r8ALC[0,5] = -540.66668
$ws.Range("H8:M8").Value = $r8ALC

$ws.Range("H29").Value = 2942.2632

$r31ALC = New-Object 'object[,]' 1,7
$r31ALC[0,0] = 525.7143
$r31ALC[0,1] = 530
$r31ALC[0,2] = 500
$r31ALC[0,3] = 1590
$r31ALC[0,4] = 1500
$r31ALC[0,5] = -1360
$r31ALC[0,6] = -1960
$ws.Range("H31:N31").Value = $r31ALC

$r43ALC = New-Object 'object[,]' 1,7
$r43ALC[0,0] = 688.4375
$r43ALC[0,1] = 406.66666
$r43ALC[0,2] = 857.5
$r43ALC[0,3] = 406.66666
$r43ALC[0,4] = 857.5
$r43ALC[0,5] = -337.66666
$r43ALC[0,6] = -995.5
$ws.Range("H43:N43").Value = $r43ALC

$r86ALC = New-Object 'object[,]' 1,7
$r86ALC[0,0] = 47638.773
$r86ALC[0,1] = 68663.53
$r86ALC[0,2] = 2585.7144
$r86ALC[0,3] = 68663.53
$r86ALC[0,4] = 2585.7144
$r86ALC[0,5] = -67540.53
$r86ALC[0,6] = -4831.7144
$ws.Range("H86:N86").Value = $r86ALC

$r88ALC = New-Object 'object[,]' 1,7
$r88ALC[0,0] = 3133.3333
$r88ALC[0,1] = 1900
$r88ALC[0,2] = 3750
$r88ALC[0,3] = 1900
$r88ALC[0,4] = 3750
$r88ALC[0,5] = -1494
$r88ALC[0,6] = -4562
$ws.Range("H88:N88").Value = $r88ALC

$r89ALC = New-Object 'object[,]' 1,7
$r89ALC[0,0] = 47638.773
$r89ALC[0,1] = 68663.53
$r89ALC[0,2] = 2585.7144
$r89ALC[0,3] = 343317.65
$r89ALC[0,4] = 12928.572
$r89ALC[0,5] = -337701.65
$r89ALC[0,6] = -24160.572
$ws.Range("H89:N89").Value = $r89ALC

$r91ALC = New-Object 'object[,]' 1,7
$r91ALC[0,0] = 3133.3333
$r91ALC[0,1] = 1900
$r91ALC[0,2] = 3750
$r91ALC[0,3] = 1900
$r91ALC[0,4] = 3750
$r91ALC[0,5] = -496
$r91ALC[0,6] = -6558
$ws.Range("H91:N91").Value = $r91ALC

$r100ALC = New-Object 'object[,]' 1,7
$r100ALC[0,0] = 2197.7273
$r100ALC[0,1] = 1655.3846
$r100ALC[0,2] = 2981.111
$r100ALC[0,3] = 1655.3846
$r100ALC[0,4] = 2981.111
$r100ALC[0,5] = -1114.3846
$r100ALC[0,6] = -4063.111
$ws.Range("H100:N100").Value = $r100ALC

$r137ALC = New-Object 'object[,]' 1,7
$r137ALC[0,0] = 9651601
$r137ALC[0,1] = 2457.5
$r137ALC[0,2] = 26192992
$r137ALC[0,3] = 7372.5
$r137ALC[0,4] = 78578976
$r137ALC[0,5] = -4822.5
$r137ALC[0,6] = -78584076
$ws.Range("H137:N137").Value = $r137ALC

$ws = $wb.Worksheets.Item("ARM")
$r11ARM = New-Object 'object[,]' 1,7
$r11ARM[0,0] = 206801.2
$r11ARM[0,1] = 1000000
$r11ARM[0,2] = 8501.5
$r11ARM[0,3] = 1000000
$r11ARM[0,4] = 8501.5
$r11ARM[0,5] = -999856
$r11ARM[0,6] = -8789.5
$ws.Range("H11:N11").Value = $r11ARM

$r88ARM = New-Object 'object[,]' 1,7
$r88ARM[0,0] = 4715.591
$r88ARM[0,1] = 7061.6
$r88ARM[0,2] = 2760.5833
$r88ARM[0,3] = 7061.6
$r88ARM[0,4] = 2760.5833
$r88ARM[0,5] = -6655.6
$r88ARM[0,6] = -3572.5833
$ws.Range("H88:N88").Value = $r88ARM

$r91ARM = New-Object 'object[,]' 1,7
$r91ARM[0,0] = 4715.591
$r91ARM[0,1] = 7061.6
$r91ARM[0,2] = 2760.5833
$r91ARM[0,3] = 7061.6
$r91ARM[0,4] = 2760.5833
$r91ARM[0,5] = -5657.6
$r91ARM[0,6] = -5568.5833
$ws.Range("H91:N91").Value = $r91ARM

$r97ARM = New-Object 'object[,]' 1,7
$r97ARM[0,0] = 1244.3334
$r97ARM[0,1] = 1154.6428
$r97ARM[0,2] = 2500
$r97ARM[0,3] = 1154.6428
$r97ARM[0,4] = 2500
$r97ARM[0,5] = -658.6428000000001
$r97ARM[0,6] = -3492
$ws.Range("H97:N97").Value = $r97ARM

$r134ARM = New-Object 'object[,]' 1,7
$r134ARM[0,0] = 0
$r134ARM[0,1] = 0
$r134ARM[0,2] = 0
$r134ARM[0,3] = 0
$r134ARM[0,4] = $null
$r134ARM[0,5] = $null
$r134ARM[0,6] = 0
$ws.Range("H134:N134").Value = $r134ARM

$r135ARM = New-Object 'object[,]' 1,7
$r135ARM[0,0] = 70000
$r135ARM[0,1] = 0
$r135ARM[0,2] = 70000
$r135ARM[0,3] = 0
$r135ARM[0,4] = 70000
$r135ARM[0,5] = $null
$r135ARM[0,6] = -80140
$ws.Range("H135:N135").Value = $r135ARM

$ws = $wb.Worksheets.Item("BSM")
$r99BSM = New-Object 'object[,]' 1,6
$r99BSM[0,0] = 2325.353
$r99BSM[0,1] = 1248.5714
$r99BSM[0,2] = 3079.1
$r99BSM[0,3] = 1248.5714
$r99BSM[0,4] = 3079.1
$r99BSM[0,5] = 249.4286
$ws.Range("H99:M99").Value = $r99BSM

$r107BSM = New-Object 'object[,]' 1,7
$r107BSM[0,0] = 1878.0454
$r107BSM[0,1] = 1629.6471
$r107BSM[0,2] = 2722.6
$r107BSM[0,3] = 1629.6471
$r107BSM[0,4] = 2722.6
$r107BSM[0,5] = 290.3529000000001
$r107BSM[0,6] = -6562.6
$ws.Range("H107:N107").Value = $r107BSM

$r135BSM = New-Object 'object[,]' 1,7
$r135BSM[0,0] = 49936
$r135BSM[0,1] = 0
$r135BSM[0,2] = 49936
$r135BSM[0,3] = 0
$r135BSM[0,4] = 49936
$r135BSM[0,5] = $null
$r135BSM[0,6] = -60076
$ws.Range("H135:N135").Value = $r135BSM

$ws = $wb.Worksheets.Item("CRP")
$r31CRP = New-Object 'object[,]' 1,7
$r31CRP[0,0] = 427974.38
$r31CRP[0,1] = 3291.7896
$r31CRP[0,2] = 757320.0600000001
$r31CRP[0,3] = 3291.7896
$r31CRP[0,4] = 757320.0600000001
$r31CRP[0,5] = -2996.7896
$r31CRP[0,6] = -757910.0600000001
$ws.Range("H31:N31").Value = $r31CRP

$r34CRP = New-Object 'object[,]' 1,7
$r34CRP[0,0] = 427974.38
$r34CRP[0,1] = 3291.7896
$r34CRP[0,2] = 757320.0600000001
$r34CRP[0,3] = 3291.7896
$r34CRP[0,4] = 757320.0600000001
$r34CRP[0,5] = -3089.7896
$r34CRP[0,6] = -757724.0600000001
$ws.Range("H34:N34").Value = $r34CRP

$r135CRP = New-Object 'object[,]' 1,7
$r135CRP[0,0] = 0
$r135CRP[0,1] = 0
$r135CRP[0,2] = 0
$r135CRP[0,3] = 0
$r135CRP[0,4] = $null
$r135CRP[0,5] = $null
$r135CRP[0,6] = 0
$ws.Range("H135:N135").Value = $r135CRP

$ws = $wb.Worksheets.Item("CUL")
$r44CUL = New-Object 'object[,]' 1,7
$r44CUL[0,0] = 497.57144
$r44CUL[0,1] = 320.75
$r44CUL[0,2] = 733.3333
$r44CUL[0,3] = 962.25
$r44CUL[0,4] = 2199.9999
$r44CUL[0,5] = -564.25
$r44CUL[0,6] = -2995.9999
$ws.Range("H44:N44").Value = $r44CUL

$r47CUL = New-Object 'object[,]' 1,7
$r47CUL[0,0] = 537.2308
$r47CUL[0,1] = 370.85715
$r47CUL[0,2] = 731.3333
$r47CUL[0,3] = 1112.57145
$r47CUL[0,4] = 2193.9999
$r47CUL[0,5] = -681.5714499999999
$r47CUL[0,6] = -3055.9999
$ws.Range("H47:N47").Value = $r47CUL

$r104CUL = New-Object 'object[,]' 1,7
$r104CUL[0,0] = 1513
$r104CUL[0,1] = 1026
$r104CUL[0,2] = 2000
$r104CUL[0,3] = 3078
$r104CUL[0,4] = 6000
$r104CUL[0,5] = -457
$r104CUL[0,6] = -11242
$ws.Range("H104:N104").Value = $r104CUL

$r129CUL = New-Object 'object[,]' 1,6
$r129CUL[0,0] = 1402.4333
$r129CUL[0,1] = 799.0909
$r129CUL[0,2] = 1751.7368
$r129CUL[0,3] = 2397.2727
$r129CUL[0,4] = 5255.2104
$r129CUL[0,5] = 2602.7273
$ws.Range("H129:M129").Value = $r129CUL

$r131CUL = New-Object 'object[,]' 1,7
$r131CUL[0,0] = 1067.2787
$r131CUL[0,1] = 225.8
$r131CUL[0,2] = 1341.674
$r131CUL[0,3] = 677.4000000000001
$r131CUL[0,4] = 4025.022
$r131CUL[0,5] = 4362.6
$r131CUL[0,6] = -14105.022
$ws.Range("H131:N131").Value = $r131CUL

$r134CUL = New-Object 'object[,]' 1,7
$r134CUL[0,0] = 3578.658
$r134CUL[0,1] = 3260.4285
$r134CUL[0,2] = 3971.7646
$r134CUL[0,3] = 9781.2855
$r134CUL[0,4] = 11915.2938
$r134CUL[0,5] = -4711.2855
$r134CUL[0,6] = -22055.2938
$ws.Range("H134:N134").Value = $r134CUL

$ws = $wb.Worksheets.Item("GSM")
$r107GSM = New-Object 'object[,]' 1,7
$r107GSM[0,0] = 706.375
$r107GSM[0,1] = 230.4
$r107GSM[0,2] = 1499.6666
$r107GSM[0,3] = 230.4
$r107GSM[0,4] = 1499.6666
$r107GSM[0,5] = 1689.6
$r107GSM[0,6] = -5339.6666
$ws.Range("H107:N107").Value = $r107GSM

$r135GSM = New-Object 'object[,]' 1,7
$r135GSM[0,0] = 0
$r135GSM[0,1] = 0
$r135GSM[0,2] = 0
$r135GSM[0,3] = 0
$r135GSM[0,4] = $null
$r135GSM[0,5] = $null
$r135GSM[0,6] = 0
$ws.Range("H135:N135").Value = $r135GSM

$ws = $wb.Worksheets.Item("LTW")
$r55LTW = New-Object 'object[,]' 1,7
$r55LTW[0,0] = 150
$r55LTW[0,1] = 150
$r55LTW[0,2] = 0
$r55LTW[0,3] = 150
$r55LTW[0,4] = 0
$r55LTW[0,5] = $null
$r55LTW[0,6] = 23
$ws.Range("H55:N55").Value = $r55LTW

$r82LTW = New-Object 'object[,]' 1,7
$r82LTW[0,0] = 1955.25
$r82LTW[0,1] = 985
$r82LTW[0,2] = 2925.5
$r82LTW[0,3] = 985
$r82LTW[0,4] = 2925.5
$r82LTW[0,5] = -624
$r82LTW[0,6] = -3647.5
$ws.Range("H82:N82").Value = $r82LTW

$r85LTW = New-Object 'object[,]' 1,7
$r85LTW[0,0] = 1955.25
$r85LTW[0,1] = 985
$r85LTW[0,2] = 2925.5
$r85LTW[0,3] = 985
$r85LTW[0,4] = 2925.5
$r85LTW[0,5] = 263
$r85LTW[0,6] = -5421.5
$ws.Range("H85:N85").Value = $r85LTW

$r134LTW = New-Object 'object[,]' 1,7
$r134LTW[0,0] = 53109
$r134LTW[0,1] = 0
$r134LTW[0,2] = 53109
$r134LTW[0,3] = 0
$r134LTW[0,4] = 53109
$r134LTW[0,5] = $null
$r134LTW[0,6] = -63249
$ws.Range("H134:N134").Value = $r134LTW

$r135LTW = New-Object 'object[,]' 1,7
$r135LTW[0,0] = 67999
$r135LTW[0,1] = 0
$r135LTW[0,2] = 67999
$r135LTW[0,3] = 0
$r135LTW[0,4] = 67999
$r135LTW[0,5] = $null
$r135LTW[0,6] = -78139
$ws.Range("H135:N135").Value = $r135LTW

$ws = $wb.Worksheets.Item("WVR")
$r62WVR = New-Object 'object[,]' 1,7
$r62WVR[0,0] = 3802.4
$r62WVR[0,1] = 0
$r62WVR[0,2] = 3802.4
$r62WVR[0,3] = 0
$r62WVR[0,4] = 3802.4
$r62WVR[0,5] = $null
$r62WVR[0,6] = -5050.4
$ws.Range("H62:N62").Value = $r62WVR

$r65WVR = New-Object 'object[,]' 1,7
$r65WVR[0,0] = 3802.4
$r65WVR[0,1] = 0
$r65WVR[0,2] = 3802.4
$r65WVR[0,3] = 0
$r65WVR[0,4] = 19012
$r65WVR[0,5] = $null
$r65WVR[0,6] = -25252
$ws.Range("H65:N65").Value = $r65WVR

$r81WVR = New-Object 'object[,]' 1,7
$r81WVR[0,0] = 11768394
$r81WVR[0,1] = 2736
$r81WVR[0,2] = 16670751
$r81WVR[0,3] = 5472
$r81WVR[0,4] = 33341502
$r81WVR[0,5] = -4411
$r81WVR[0,6] = -33343624
$ws.Range("H81:N81").Value = $r81WVR

$r84WVR = New-Object 'object[,]' 1,7
$r84WVR[0,0] = 11768394
$r84WVR[0,1] = 2736
$r84WVR[0,2] = 16670751
$r84WVR[0,3] = 27360
$r84WVR[0,4] = 166707510
$r84WVR[0,5] = -22056
$r84WVR[0,6] = -166718118
$ws.Range("H84:N84").Value = $r84WVR

$r132WVR = New-Object 'object[,]' 1,7
$r132WVR[0,0] = 2894.1177
$r132WVR[0,1] = 3116.75
$r132WVR[0,2] = 2359.8
$r132WVR[0,3] = 9350.25
$r132WVR[0,4] = 7079.400000000001
$r132WVR[0,5] = -6820.25
$r132WVR[0,6] = -12139.4
$ws.Range("H132:N132").Value = $r132WVR

$r135WVR = New-Object 'object[,]' 1,7
$r135WVR[0,0] = 41000
$r135WVR[0,1] = 0
$r135WVR[0,2] = 41000
$r135WVR[0,3] = 0
$r135WVR[0,4] = 41000
$r135WVR[0,5] = $null
$r135WVR[0,6] = -51140
$ws.Range("H135:N135").Value = $r135WVR
